$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1983805668016194
$ws.Cells.Item(2, 3).Value = 0.5627530364372469
$ws.Cells.Item(2, 10).Value = 0.01619433198380567
$ws.Cells.Item(2, 16).Value = 0.1376518218623482
$ws.Cells.Item(2, 19).Value = 0.08502024291497975
$ws.Cells.Item(3, 2).Value = 0.01408450704225352
$ws.Cells.Item(3, 10).Value = 0.02816901408450704
$ws.Cells.Item(3, 16).Value = 0.7887323943661971
$ws.Cells.Item(3, 19).Value = 0.1690140845070423
$ws.Cells.Item(4, 10).Value = 0.08108108108108109
$ws.Cells.Item(4, 16).Value = 0.5945945945945946
$ws.Cells.Item(4, 19).Value = 0.3243243243243243
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(6, 2).Value = 0.03255813953488372
$ws.Cells.Item(6, 4).Value = 0.0186046511627907
$ws.Cells.Item(6, 6).Value = 0.06976744186046512
$ws.Cells.Item(6, 10).Value = 0.3395348837209302
$ws.Cells.Item(6, 15).Value = 0.0186046511627907
$ws.Cells.Item(6, 17).Value = 0.1441860465116279
$ws.Cells.Item(6, 18).Value = 0.06511627906976744
$ws.Cells.Item(6, 19).Value = 0.3116279069767442
$ws.Cells.Item(7, 2).Value = 0.06179775280898876
$ws.Cells.Item(7, 4).Value = 0.02247191011235955
$ws.Cells.Item(7, 5).Value = 0.005617977528089887
$ws.Cells.Item(7, 6).Value = 0.05617977528089887
$ws.Cells.Item(7, 10).Value = 0.1573033707865168
$ws.Cells.Item(7, 15).Value = 0.01685393258426966
$ws.Cells.Item(7, 17).Value = 0.1292134831460674
$ws.Cells.Item(7, 18).Value = 0.1348314606741573
$ws.Cells.Item(7, 19).Value = 0.4157303370786517
$ws.Cells.Item(8, 2).Value = 0.07235142118863049
$ws.Cells.Item(8, 4).Value = 0.01808785529715762
$ws.Cells.Item(8, 5).Value = 0.002583979328165375
$ws.Cells.Item(8, 6).Value = 0.041343669250646
$ws.Cells.Item(8, 10).Value = 0.1291989664082687
$ws.Cells.Item(8, 15).Value = 0.0103359173126615
$ws.Cells.Item(8, 17).Value = 0.1679586563307494
$ws.Cells.Item(8, 18).Value = 0.1266149870801034
$ws.Cells.Item(8, 19).Value = 0.4315245478036176
$ws.Cells.Item(9, 2).Value = 0.1117021276595745
$ws.Cells.Item(9, 4).Value = 0.01595744680851064
$ws.Cells.Item(9, 6).Value = 0.06914893617021277
$ws.Cells.Item(9, 10).Value = 0.1063829787234043
$ws.Cells.Item(9, 15).Value = 0.02127659574468085
$ws.Cells.Item(9, 17).Value = 0.1595744680851064
$ws.Cells.Item(9, 18).Value = 0.101063829787234
$ws.Cells.Item(9, 19).Value = 0.4148936170212766
$ws.Cells.Item(10, 2).Value = 0.1078348778433024
$ws.Cells.Item(10, 4).Value = 0.0160067396798652
$ws.Cells.Item(10, 6).Value = 0.07497893850042123
$ws.Cells.Item(10, 10).Value = 0.1204717775905644
$ws.Cells.Item(10, 15).Value = 0.02358887952822241
$ws.Cells.Item(10, 17).Value = 0.1946082561078349
$ws.Cells.Item(10, 18).Value = 0.08845829823083404
$ws.Cells.Item(10, 19).Value = 0.3740522325189554
$ws.Cells.Item(11, 7).Value = 0.1360294117647059
$ws.Cells.Item(11, 10).Value = 0.07352941176470588
$ws.Cells.Item(11, 11).Value = 0.2095588235294118
$ws.Cells.Item(11, 12).Value = 0.5588235294117647
$ws.Cells.Item(11, 19).Value = 0.02205882352941177
$ws.Cells.Item(12, 7).Value = 0.7639751552795031
$ws.Cells.Item(12, 10).Value = 0.1677018633540373
$ws.Cells.Item(12, 11).Value = 0.0124223602484472
$ws.Cells.Item(12, 12).Value = 0.03726708074534162
$ws.Cells.Item(12, 19).Value = 0.01863354037267081
$ws.Cells.Item(13, 7).Value = 0.696969696969697
$ws.Cells.Item(13, 10).Value = 0.2424242424242424
$ws.Cells.Item(13, 19).Value = 0.06060606060606061
$ws.Cells.Item(15, 6).Value = 0.04326923076923077
$ws.Cells.Item(15, 8).Value = 0.1538461538461539
$ws.Cells.Item(15, 9).Value = 0.07211538461538461
$ws.Cells.Item(15, 10).Value = 0.3557692307692308
$ws.Cells.Item(15, 11).Value = 0.05288461538461538
$ws.Cells.Item(15, 13).Value = 0.004807692307692308
$ws.Cells.Item(15, 14).Value = 0.004807692307692308
$ws.Cells.Item(15, 15).Value = 0.0625
$ws.Cells.Item(15, 19).Value = 0.25
$ws.Cells.Item(16, 6).Value = 0.01886792452830189
$ws.Cells.Item(16, 8).Value = 0.1886792452830189
$ws.Cells.Item(16, 9).Value = 0.09433962264150944
$ws.Cells.Item(16, 10).Value = 0.3710691823899371
$ws.Cells.Item(16, 11).Value = 0.1069182389937107
$ws.Cells.Item(16, 13).Value = 0.01257861635220126
$ws.Cells.Item(16, 15).Value = 0.03773584905660377
$ws.Cells.Item(16, 19).Value = 0.169811320754717
$ws.Cells.Item(17, 6).Value = 0.02362204724409449
$ws.Cells.Item(17, 8).Value = 0.1758530183727034
$ws.Cells.Item(17, 9).Value = 0.1076115485564304
$ws.Cells.Item(17, 10).Value = 0.4041994750656168
$ws.Cells.Item(17, 11).Value = 0.08398950131233596
$ws.Cells.Item(17, 13).Value = 0.01837270341207349
$ws.Cells.Item(17, 15).Value = 0.05774278215223097
$ws.Cells.Item(17, 19).Value = 0.1286089238845144
$ws.Cells.Item(18, 6).Value = 0.04265402843601896
$ws.Cells.Item(18, 8).Value = 0.1658767772511848
$ws.Cells.Item(18, 9).Value = 0.0947867298578199
$ws.Cells.Item(18, 10).Value = 0.3601895734597156
$ws.Cells.Item(18, 11).Value = 0.0947867298578199
$ws.Cells.Item(18, 13).Value = 0.01895734597156398
$ws.Cells.Item(18, 15).Value = 0.07582938388625593
$ws.Cells.Item(18, 19).Value = 0.1469194312796208
$ws.Cells.Item(19, 6).Value = 0.01524132091447926
$ws.Cells.Item(19, 8).Value = 0.1947502116850127
$ws.Cells.Item(19, 9).Value = 0.0821337849280271
$ws.Cells.Item(19, 10).Value = 0.388653683319221
$ws.Cells.Item(19, 11).Value = 0.1117696867061812
$ws.Cells.Item(19, 13).Value = 0.01947502116850127
$ws.Cells.Item(19, 15).Value = 0.07197290431837426
$ws.Cells.Item(19, 19).Value = 0.1160033869602032

Write-Host "Applied 109 cell updates"